# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 304 in the "Ajo" (garlic) sheet,
# shifting all subsequent rows down by one (dimension grows from R323 to R324).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 304.
$ws.Rows.Item(304).Insert()

# Populate the new row 304 with the new weekly record.
$ws.Range("A304").Value = 7
$ws.Range("B304").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C304").Value = "Ñuble"
$ws.Range("D304").Value = 44931
$ws.Range("E304").Value = 16
$ws.Range("F304").Value = 100112003
$ws.Range("G304").Value = "Ajo"
$ws.Range("H304").Value = "Chino"
$ws.Range("I304").Value = "Primera"
$ws.Range("J304").Value = 60
$ws.Range("K304").Value = 15000
$ws.Range("L304").Value = 16000
$ws.Range("M304").Value = 15500
$ws.Range("N304").Value = "$/caja 10 kilos"
$ws.Range("O304").Value = "China"
$ws.Range("P304").Value = 1550
$ws.Range("Q304").Value = 10
$ws.Range("R304").Value = "Hortaliza"
